$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.739.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.07%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.285.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.48%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'97.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.20%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'269.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.61%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.29%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -2.57%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'45.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.34%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.22%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -2.95%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +2.37%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'15.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.33%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.627.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.53%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.855"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.12%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.296.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.45%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'43.775.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.20%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +2.35%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.36%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'72.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.31%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'2.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +8.77%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'232.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.34%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'9.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.44%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +5.40%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.04%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'11.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.51%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'3.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.74%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.32%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'38.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.84%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'175.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.13%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'21.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.25%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0895"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.11%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.72%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.48%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +7.07%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.40%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0351"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.18%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.70%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.19%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.91%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'12.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.59%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.12%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'64.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.17%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -3.23%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'8.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.23%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.81%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Aave"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'98.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.87%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'TrustWalletToken"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'1.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.13%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.440"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.58%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +10.37%  "
$ws.Range("E51").Style = "Normal"
